$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.113.11'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.592.15'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("D4").Value = "'0.9967"
$ws.Range("E4").Value = '  -0.36%  '
$ws.Range("D5").Value = "'0.9972"
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = "'0.3762"
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("D8").Value = "'0.3612"
$ws.Range("E8").Value = '  -1.23%  '
$ws.Range("D9").Value = "'51.10"
$ws.Range("E9").Value = '  +3.33%  '
$ws.Range("D10").Value = "'1.247"
$ws.Range("E10").Value = '  -1.57%  '
$ws.Range("D11").Value = "'0.9968"
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("D12").Value = "'0.08046"
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("D13").Value = "'22.42"
$ws.Range("E13").Value = '  -2.45%  '
$ws.Range("D14").Value = "'6.524"
$ws.Range("E14").Value = '  -1.19%  '
$ws.Range("D15").Value = "'7.364"
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("D16").Value = "'0.00001239"
$ws.Range("E16").Value = '  -1.13%  '
$ws.Range("D17").Value = '1.594.43'
$ws.Range("E17").Value = '  -1.01%  '
$ws.Range("D18").Value = "'92.93"
$ws.Range("E18").Value = '  +1.67%  '
$ws.Range("D19").Value = "'0.06767"
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").Value = "'17.93"
$ws.Range("E20").Value = '  -2.35%  '
$ws.Range("D21").Value = "'6.448"
$ws.Range("E21").Value = '  -1.89%  '
$ws.Range("D22").Value = "'0.9974"
$ws.Range("E22").Value = '  -0.45%  '
$ws.Range("D23").Value = "'12.81"
$ws.Range("E23").Value = '  -1.87%  '
$ws.Range("D24").Value = '23.093.58'
$ws.Range("E24").Value = '  -0.36%  '
$ws.Range("D25").Value = "'2.383"
$ws.Range("E25").Value = '  +1.36%  '
$ws.Range("D26").Value = "'2.916"
$ws.Range("E26").Value = '  +2.50%  '
$ws.Range("D27").Value = "'20.90"
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("D28").Value = "'148.54"
$ws.Range("E28").Value = '  -1.32%  '
$ws.Range("D29").Value = "'5.186"
$ws.Range("E29").Value = '  -1.77%  '
$ws.Range("D30").Value = "'133.18"
$ws.Range("E30").Value = '  +0.36%  '
$ws.Range("D31").Value = "'2.372"
$ws.Range("E31").Value = '  -1.55%  '
$ws.Range("D32").Value = "'6.752"
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("D33").Value = '1.769.04'
$ws.Range("E33").Value = '  -1.12%  '
$ws.Range("D34").Value = "'0.9572"
$ws.Range("E34").Value = '  -1.34%  '
$ws.Range("D35").Value = "'0.07471"
$ws.Range("E35").Value = '  -2.88%  '
$ws.Range("D36").Value = "'10.11"
$ws.Range("E36").Value = '  -0.57%  '
$ws.Range("D37").Value = "'0.02676"
$ws.Range("E37").Value = '  -3.18%  '
$ws.Range("D38").Value = "'0.2500"
$ws.Range("E38").Value = '  -2.38%  '
$ws.Range("D39").Value = "'6.101"
$ws.Range("E39").Value = '  -2.23%  '
$ws.Range("D40").Value = "'0.08784"
$ws.Range("E40").Value = '  -1.32%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = "'0.7094"
$ws.Range("E41").Value = '  -1.38%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = "'1.356"
$ws.Range("E42").Value = '  -2.55%  '
$ws.Range("D43").Value = "'12.20"
$ws.Range("E43").Value = '  -4.85%  '
$ws.Range("D44").Value = "'14.95"
$ws.Range("E44").Value = '  -4.87%  '
$ws.Range("D45").Value = "'0.6513"
$ws.Range("E45").Value = '  -2.33%  '
$ws.Range("D46").Value = "'0.9962"
$ws.Range("E46").Value = '  -0.31%  '
$ws.Range("D47").Value = "'3.993"
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("D48").Value = "'2.283"
$ws.Range("E48").Value = '  -1.04%  '
$ws.Range("D49").Value = "'131.43"
$ws.Range("E49").Value = '  +0.53%  '
$ws.Range("D50").Value = "'0.07903"
$ws.Range("E50").Value = '  -1.28%  '
$ws.Range("D51").Value = "'1.222"
$ws.Range("E51").Value = '  +3.76%  '

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
